# Update countries & provincias Spain
# Updates the COVID country-stats sheet with refreshed numbers for the
# 23-Oct-2020 21:04 data pull, and re-applies the (sort-induced) row swaps
# for three pairs of adjacent countries whose totals crossed rank:
#   Corea del Sur / Republica de Macedonia
#   Camboya / Monaco
#   Montserrat / Islas Malvinas

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 21:04"

# Row 4 - Estados Unidos (new totals)
$ws.Range("B4").Value = 8699988
$ws.Range("C4").Value = 38337
$ws.Range("D4").Value = 5672924
$ws.Range("E4").Value = 2798250
$ws.Range("G4").Value = 433
$ws.Range("H4").Value = 228814

# Row 5 - India (new totals)
$ws.Range("B5").Value = 7813667
$ws.Range("C5").Value = 54027
$ws.Range("D5").Value = 7013569
$ws.Range("E5").Value = 682106
$ws.Range("G5").Value = 656
$ws.Range("H5").Value = 117992

# Row 10 - Francia (new totals)
$ws.Range("B10").Value = 1041075
$ws.Range("C10").Value = 42032
$ws.Range("D10").Value = 109486
$ws.Range("E10").Value = 897081
$ws.Range("G10").Value = 298
$ws.Range("H10").Value = 34508

# Row 20 - Alemania (new totals)
$ws.Range("B20").Value = 416410
$ws.Range("C20").Value = 12536
$ws.Range("E20").Value = 96126
$ws.Range("G20").Value = 40
$ws.Range("H20").Value = 10084

# Row 24 - Turquia (new totals)
$ws.Range("B24").Value = 357693
$ws.Range("C24").Value = 2165
$ws.Range("D24").Value = 311520
$ws.Range("E24").Value = 36515
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = 9658

# Row 50 - Suiza (new totals)
$ws.Range("E50").Value = 45786
$ws.Range("G50").Value = 15
$ws.Range("H50").Value = 2067

# Rows 89/90 - Corea del Sur and Republica de Macedonia swap rank:
# Republica de Macedonia moves up to row 89 with fresh numbers, Corea del
# Sur drops to row 90 carrying its unchanged prior totals.
$ws.Range("A89").Value = "Republica de Macedonia"
$ws.Range("B89").Value = 25991
$ws.Range("C89").Value = 518
$ws.Range("D89").Value = 18247
$ws.Range("E89").Value = 6861
$ws.Range("G89").Value = 9
$ws.Range("H89").Value = 883

$ws.Range("A90").Value = "Corea del Sur"
$ws.Range("B90").Value = 25698
$ws.Range("C90").Value = 155
$ws.Range("D90").Value = 23717
$ws.Range("E90").Value = 1526
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = 455

# Row 128 - Republica de Yibuti (new totals)
$ws.Range("B128").Value = 5528
$ws.Range("C128").Value = 6
$ws.Range("D128").Value = 5393
$ws.Range("E128").Value = 74

# Row 165 - Republica del Chad (new totals)
$ws.Range("B165").Value = 1423
$ws.Range("C165").Value = 13
$ws.Range("D165").Value = 1234
$ws.Range("E165").Value = 93

# Row 166 - Liberia (new totals)
$ws.Range("B166").Value = 1393
$ws.Range("C166").Value = 8
$ws.Range("E166").Value = 33

# Rows 189/190 - Camboya and Monaco swap rank:
# Monaco moves up to row 189 with fresh numbers, Camboya drops to row 190
# carrying its unchanged prior totals.
$ws.Range("A189").Value = "Monaco"
$ws.Range("B189").Value = 295
$ws.Range("C189").Value = 14
$ws.Range("D189").Value = 241
$ws.Range("E189").Value = 52
$ws.Range("H189").Value = 2

$ws.Range("A190").Value = "Camboya"
$ws.Range("B190").Value = 286
$ws.Range("D190").Value = 280
$ws.Range("E190").Value = 6
$ws.Range("H190").Value = 0

# Rows 216/217 - Montserrat and Islas Malvinas swap rank (totals unchanged,
# just traded places).
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0

$ws.Range("A217").Value = "Montserrat"
$ws.Range("D217").Value = 12
$ws.Range("H217").Value = 1
